$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number for rows 2 through 24.
# Update every one of those cells from 45222 (2023-10-23) to 45224 (2023-10-25).
for ($row = 2; $row -le 24; $row++) {
    $ws.Cells.Item($row, 3).Value = 45224
}
